$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing column (L) into the new column (M)
# for the border/thick-bottom rows so the new column matches the table style.
$ws.Range("L3").Copy() | Out-Null
$ws.Range("M3").PasteSpecial(-4122) | Out-Null

$ws.Range("L4").Copy() | Out-Null
$ws.Range("M4").PasteSpecial(-4122) | Out-Null

$ws.Range("L5").Copy() | Out-Null
$ws.Range("M5").PasteSpecial(-4122) | Out-Null

# New data values for the added 2022 column
$ws.Range("M4").Value = 2022
$ws.Range("M5").Value = 373

# Update the active selection as recorded in the sheet view
$ws.Range("O4").Select() | Out-Null
